# Auto-generated edit script applying scheduled-runner updates to market-price
# statistic columns (H/I/J/K/L/M/N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1732.619
$ws.Range("J17").Value = 2580.4546
$ws.Range("L17").Value = 7741.3638
$ws.Range("N17").Value = -8077.3638
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
# Row 92
$ws.Range("H92").Value = 168
$ws.Range("I92").Value = 174.9
$ws.Range("K92").Value = 174.9
$ws.Range("M92").Value = 1073.1
# Row 116
$ws.Range("H116").Value = 17860.455
$ws.Range("I116").Value = 17209.285
$ws.Range("J116").Value = 19000
$ws.Range("K116").Value = 17209.285
$ws.Range("L116").Value = 19000
$ws.Range("M116").Value = -13767.285
$ws.Range("N116").Value = -25884
# Row 129
$ws.Range("H129").Value = 889.94446
$ws.Range("I129").Value = 453.4
$ws.Range("K129").Value = 1360.2
$ws.Range("M129").Value = 3639.8
# Row 135
$ws.Range("H135").Value = 664.5333000000001
$ws.Range("I135").Value = 564.913
$ws.Range("K135").Value = 5084.217000000001
$ws.Range("M135").Value = -2549.217000000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 799.06665
$ws.Range("I2").Value = 799.06665
$ws.Range("K2").Value = 799.06665
$ws.Range("M2").Value = -686.06665
# Row 45
$ws.Range("H45").Value = 2732.889
$ws.Range("I45").Value = 2073.6667
$ws.Range("K45").Value = 2073.6667
$ws.Range("M45").Value = -1696.6667
# Row 74
$ws.Range("H74").Value = 2952.2727
$ws.Range("I74").Value = 2571.5789
$ws.Range("K74").Value = 2571.5789
$ws.Range("M74").Value = -1697.5789
# Row 77
$ws.Range("H77").Value = 2952.2727
$ws.Range("I77").Value = 2571.5789
$ws.Range("K77").Value = 12857.8945
$ws.Range("M77").Value = -8489.8945
# Row 97
$ws.Range("H97").Value = 1462.4667
$ws.Range("I97").Value = 1412.2727
$ws.Range("K97").Value = 1412.2727
$ws.Range("M97").Value = -916.2727
# Row 110
$ws.Range("H110").Value = 76925576
$ws.Range("I110").Value = 125001950
$ws.Range("K110").Value = 125001950
$ws.Range("M110").Value = -124999905
# Row 116
$ws.Range("H116").Value = 799.06665
$ws.Range("I116").Value = 799.06665
$ws.Range("K116").Value = 799.06665
$ws.Range("M116").Value = 1494.93335
# Row 132
$ws.Range("H132").Value = 3553.9473
$ws.Range("I132").Value = 3473.611
$ws.Range("K132").Value = 10420.833
$ws.Range("M132").Value = -7890.832999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 799.06665
$ws.Range("I3").Value = 799.06665
$ws.Range("K3").Value = 799.06665
$ws.Range("M3").Value = -685.06665
# Row 86
$ws.Range("H86").Value = 2679.8
$ws.Range("I86").Value = 2674.75
$ws.Range("K86").Value = 2674.75
$ws.Range("M86").Value = -1551.75
# Row 89
$ws.Range("H89").Value = 2679.8
$ws.Range("I89").Value = 2674.75
$ws.Range("K89").Value = 13373.75
$ws.Range("M89").Value = -7757.75
# Row 97
$ws.Range("H97").Value = 4825.4
$ws.Range("I97").Value = 4825.4
$ws.Range("K97").Value = 4825.4
$ws.Range("M97").Value = -3834.4
# Row 99
$ws.Range("H99").Value = 125001090
$ws.Range("I99").Value = 166667760
$ws.Range("K99").Value = 166667760
$ws.Range("M99").Value = -166666262
# Row 105
$ws.Range("H105").Value = 12988334
$ws.Range("I105").Value = 15152806
$ws.Range("K105").Value = 15152806
$ws.Range("M105").Value = -15151059
# Row 134
$ws.Range("H134").Value = 2727.5676
$ws.Range("I134").Value = 845.9355
$ws.Range("J134").Value = 12449.333
$ws.Range("K134").Value = 2537.8065
$ws.Range("L134").Value = 37347.999
$ws.Range("M134").Value = -2.806500000000142
$ws.Range("N134").Value = -42417.999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 99.454544
$ws.Range("I7").Value = 70.71429000000001
$ws.Range("K7").Value = 70.71429000000001
$ws.Range("M7").Value = 42.28570999999999
# Row 16
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 400
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -113
$ws.Range("N16").ClearContents()
# Row 31
$ws.Range("H31").Value = 6090.2583
$ws.Range("I31").Value = 2799.6667
$ws.Range("J31").Value = 6880
$ws.Range("K31").Value = 2799.6667
$ws.Range("L31").Value = 6880
$ws.Range("M31").Value = -2504.6667
$ws.Range("N31").Value = -7470
# Row 34
$ws.Range("H34").Value = 6090.2583
$ws.Range("I34").Value = 2799.6667
$ws.Range("J34").Value = 6880
$ws.Range("K34").Value = 2799.6667
$ws.Range("L34").Value = 6880
$ws.Range("M34").Value = -2597.6667
$ws.Range("N34").Value = -7284
# Row 113
$ws.Range("H113").Value = 400
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 400
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1770
$ws.Range("N113").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 396.85715
$ws.Range("I12").Value = 499
$ws.Range("J12").Value = 379.83334
$ws.Range("K12").Value = 1497
$ws.Range("L12").Value = 1139.50002
$ws.Range("M12").Value = -1324
$ws.Range("N12").Value = -1485.50002
# Row 23
$ws.Range("H23").Value = 317.5
$ws.Range("J23").Value = 163.75
$ws.Range("L23").Value = 491.25
$ws.Range("N23").Value = -961.25
# Row 92
$ws.Range("H92").Value = 4699.8335
$ws.Range("I92").Value = 1833.3334
$ws.Range("K92").Value = 5500.0002
$ws.Range("M92").Value = -4252.0002
# Row 98
$ws.Range("H98").Value = 316.83334
$ws.Range("I98").Value = 300.5
$ws.Range("J98").Value = 333.16666
$ws.Range("K98").Value = 901.5
$ws.Range("L98").Value = 999.4999799999999
$ws.Range("M98").Value = 596.5
$ws.Range("N98").Value = -3995.49998
# Row 116
$ws.Range("H116").Value = 1826.3334
$ws.Range("I116").Value = 1826.3334
$ws.Range("K116").Value = 5479.0002
$ws.Range("M116").Value = -2037.0002
# Row 131
$ws.Range("H131").Value = 2360.3076
$ws.Range("I131").Value = 1959.8334
$ws.Range("J131").Value = 2703.5715
$ws.Range("K131").Value = 5879.5002
$ws.Range("L131").Value = 8110.7145
$ws.Range("M131").Value = -839.5002000000004
$ws.Range("N131").Value = -18190.7145

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 4566.6665
$ws.Range("I43").Value = 4566.6665
$ws.Range("K43").Value = 4566.6665
$ws.Range("M43").Value = -4415.6665
# Row 80
$ws.Range("H80").Value = 3906
$ws.Range("I80").Value = 3892.5715
$ws.Range("K80").Value = 3892.5715
$ws.Range("M80").Value = -2894.5715
# Row 83
$ws.Range("H83").Value = 3906
$ws.Range("I83").Value = 3892.5715
$ws.Range("K83").Value = 19462.8575
$ws.Range("M83").Value = -14470.8575
# Row 122
$ws.Range("H122").Value = 1951.5
$ws.Range("I122").Value = 1070.1428
$ws.Range("K122").Value = 3210.4284
$ws.Range("M122").Value = -760.4284000000002

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3546.4
$ws.Range("I46").Value = 797
$ws.Range("J46").Value = 4233.75
$ws.Range("K46").Value = 797
$ws.Range("L46").Value = 4233.75
$ws.Range("M46").Value = -609
$ws.Range("N46").Value = -4609.75
# Row 93
$ws.Range("H93").Value = 2283
$ws.Range("I93").Value = 2283
$ws.Range("K93").Value = 2283
$ws.Range("M93").Value = -1035

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 107
$ws.Range("H107").Value = 37037930
$ws.Range("I107").Value = 47619904
$ws.Range("K107").Value = 142859712
$ws.Range("M107").Value = -142857792
# Row 122
$ws.Range("H122").Value = 2300.3845
$ws.Range("I122").Value = 1369.375
$ws.Range("J122").Value = 3790
$ws.Range("K122").Value = 4108.125
$ws.Range("L122").Value = 11370
$ws.Range("M122").Value = -1658.125
$ws.Range("N122").Value = -16270
# Row 132
$ws.Range("H132").Value = 1808.0476
$ws.Range("I132").Value = 1351.1765
$ws.Range("K132").Value = 4053.5295
$ws.Range("M132").Value = -1523.5295
# Row 136
$ws.Range("H136").Value = 3590.3572
$ws.Range("I136").Value = 2560.9412
$ws.Range("K136").Value = 7682.823600000001
$ws.Range("M136").Value = -5132.823600000001
